# Reorders the "Periodo Mora" (column E) values for rows 16-22 so that the
# list of periods 2203..2209 appears in descending order (2209 first,
# 2203 last), and keeps the "Valor Mora" (column F) figures attached to the
# correct period by swapping the F16/F22 values accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New period order for rows 16..22 (was 2203,2204,2205,2206,2207,2208,2209)
$ws.Cells.Item(16, 5).Value = "2209"
$ws.Cells.Item(17, 5).Value = "2208"
$ws.Cells.Item(18, 5).Value = "2207"
$ws.Cells.Item(19, 5).Value = "2206"
$ws.Cells.Item(20, 5).Value = "2205"
$ws.Cells.Item(21, 5).Value = "2204"
$ws.Cells.Item(22, 5).Value = "2203"

# Valor Mora values follow their period: the 2209 row now carries 34666,
# the 2203 row now carries 40000 (the other periods keep 40000).
$ws.Cells.Item(16, 6).Value = 34666
$ws.Cells.Item(22, 6).Value = 40000
